$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: unit cost change (CONN Header 2POS 0.1 GOLD PCB) ---
$ws.Range("K4").Value = 0.29092500000000004

# --- Row 6: remove the stray "208 v00, " note in J6 ---
$ws.Range("J6").ClearContents()

# --- Insert two new rows at 8:9, pushing the Keystone CONN row (old row 8) down to
#     row 10 and the TOTAL row (old row 10) down to row 12 ---
$ws.Rows("8:9").Insert()

# --- New row 8: Thumb Screw, 6-32 x 12mm, Royal Blue ---
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 388
$ws.Range("E8").Value = "Thumb Screw, 6-32 x 12mm, Royal Blue"
$ws.Range("I8").Value = "AliExpress"
$ws.Range("J8").Value = "https://www.aliexpress.us/item/3256803687287578.html"
$ws.Range("K8").Value = 1.3133333333333332
$ws.Range("L8").Formula = "=B8*K8"

# --- New row 9: CONN PC PIN CIRC 0.060 DIA GOLD (Mill-Max), the new primary part ---
$ws.Range("B9").Value = 42
$ws.Range("C9").Value = 389
$ws.Range("E9").Value = "CONN PC PIN CIRC 0.060 DIA GOLD"
$ws.Range("G9").Value = "Mill-Max Manufacturing Inc."
$ws.Range("H9").Value = "4357-0-00-15-00-00-03-0"
$ws.Range("I9").Value = "Digikey"
$ws.Range("J9").Value = "https://www.digikey.com/en/products/detail/mill-max-manufacturing-corp/4357-0-00-15-00-00-03-0/5820224"
$ws.Range("K9").Value = 0.38764356349206347
$ws.Range("L9").Formula = "=B9*K9"

# --- Row 10 (was row 8 before the insert): Keystone CONN PC PIN becomes the alternate
#     part - qty drops to 0, the UPN/part-ver columns are cleared, and unit cost updates ---
$ws.Range("B10").Value = 0
$ws.Range("C10").ClearContents()
$ws.Range("F10").ClearContents()
$ws.Range("K10").Value = 0.92167827380952383

# --- Row 11: note explaining the Keystone/Mill-Max alternate relationship ---
$ws.Range("G11").Value = "Keystone is an alternate for Mill-Max 4357-0-00-15-00-00-03-0"

# --- Row 12 (was row 10 before the insert): TOTAL row - restrict the sum back to the
#     primary parts only (exclude the alternate-part row) ---
$ws.Range("L12").Formula = "=SUM(L2:L9)"

# --- Selection, matching the saved workbook view ---
$ws.Range("J6").Select()
